$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-assert the pre-existing heading formatting (large title cell / bold
# section heading) so it keeps rendering the way it already did before
# this edit - we are only adding the size-class table below, nothing
# about the existing Suriname / MSME Participation heading should change.
$ws.Range("A1").Font.Size = 18
$ws.Range("A3").Font.Bold = $true

# --- MSME size-class breakdown table (rows 9-13, cols A-D) ---

# Header row (row 9): column labels for the three metrics, bold like the
# existing "title" style already used for the section heading in A3.
$ws.Range("B9").Value = "Number of employees"
$ws.Range("C9").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D9").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B9:D9").Font.Bold = $true

# Micro (row 10)
$ws.Range("A10").Value = "Micro"
$ws.Range("B10").Value = "1-10"

# Small (row 11)
$ws.Range("A11").Value = "Small"
$ws.Range("B11").Value = "11-25"

# Medium (row 12)
$ws.Range("A12").Value = "Medium"
$ws.Range("B12").Value = "26-100"

# Large (row 13)
$ws.Range("A13").Value = "Large"
$ws.Range("B13").Value = ">100"
